# Generate Report for Handback
# Updates the existing handback row (new GUID-named source file that replaced
# the previous run's file) and appends a brand-new row for a second file that
# was handed back in the same batch, across all three sheets (Overview,
# zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$oldGuid = "3a23fbe1-3276-4940-9189-2563b186e97d"
$guid1   = "8c4fc406-04cd-4e0f-adc6-88c7333da1f4"
$guid2   = "8fd3d8bc-dd35-4bc2-bb39-61b7e6ce75db"

$hash1zh = "6141e882d664c8cb279a8c3fbedf4d4b4198bb00"
$hash2zh = "673b5204cdf159f930489f7bb909aac0c8d75e19"

# ----------------------------------------------------------------------
# Overview sheet
# ----------------------------------------------------------------------
$wsOv = $wb.Worksheets.Item("Overview")
$wsOv.Hyperlinks.Delete()

$loOv = $wsOv.ListObjects.Item(1)
$rowOv3 = $loOv.ListRows.Add()

# Row 2 (existing file re-run with a new GUID name)
$wsOv.Range("A2").Value = "$guid1.md"
$wsOv.Range("B2").Value = "e2e\$guid1.md"
$wsOv.Range("C2").Value = ".md"
$wsOv.Range("E2").Value = "Handed back: in sync with en-US"
$wsOv.Range("F2").Value = "Handed back: in sync with en-US"
$wsOv.Range("G2").Value = "2016-08-12 17:16:36"

# Row 3 (new file)
$wsOv.Range("A3").Value = "$guid2.md"
$wsOv.Range("B3").Value = "e2e\$guid2.md"
$wsOv.Range("C3").Value = ".md"
$wsOv.Range("E3").Value = "Handed back: in sync with en-US"
$wsOv.Range("F3").Value = "Handed back: in sync with en-US"
$wsOv.Range("G3").Value = "2016-08-12 17:16:36"

$wsOv.Range("G2:G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOv.Hyperlinks.Add($wsOv.Range("B2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/6fa018cb586824fa5ee6519a285314c1c36169a3/e2e/$guid1.md", "", "", "e2e\$guid1.md") | Out-Null
$wsOv.Hyperlinks.Add($wsOv.Range("B3"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/6fa018cb586824fa5ee6519a285314c1c36169a3/e2e/$guid2.md", "", "", "e2e\$guid2.md") | Out-Null

# ----------------------------------------------------------------------
# zh-cn sheet
# ----------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Hyperlinks.Delete()

$loZh = $wsZh.ListObjects.Item(1)
$rowZh3 = $loZh.ListRows.Add()

# Row 2 (existing file re-run with a new GUID name)
$wsZh.Range("A2").Value = "$guid1.md"
$wsZh.Range("B2").Value = ".md"
$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsZh.Range("D2").Value = "e2e"
$wsZh.Range("E2").Value = "ht"
$wsZh.Range("F2").Value = "False"
$wsZh.Range("G2").Value = "$guid1.$hash1zh.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-08-12 17:16:29"
$wsZh.Range("I2").Value = "$guid1.md"
$wsZh.Range("J2").Value = "$guid1.$hash1zh.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-12 17:17:12"
$wsZh.Range("L2").Value = ""
$wsZh.Range("M2").Value = "True"
$wsZh.Range("N2").Value = ""
$wsZh.Range("O2").Value = "False"
$wsZh.Range("P2").Value = ""

# Row 3 (new file)
$wsZh.Range("A3").Value = "$guid2.md"
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "True"
$wsZh.Range("G3").Value = "$guid2.$hash2zh.zh-cn.xlf"
$wsZh.Range("H3").Value = "2016-08-12 17:16:29"
$wsZh.Range("I3").Value = "$guid2.md"
$wsZh.Range("J3").Value = "$guid2.$hash2zh.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-12 17:17:12"
$wsZh.Range("L3").Value = ""
$wsZh.Range("M3").Value = "True"
$wsZh.Range("N3").Value = ""
$wsZh.Range("O3").Value = "False"
$wsZh.Range("P3").Value = ""

$wsZh.Range("H2:H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("K2:K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/6fa018cb586824fa5ee6519a285314c1c36169a3/e2e/$guid1.md", "", "", "$guid1.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/ec391cc42ee51fcad19cc5d0a141975e9ece5477/e2e/$guid1.md", "", "", "$guid1.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/6fa018cb586824fa5ee6519a285314c1c36169a3/e2e/$guid2.md", "", "", "$guid2.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/ec391cc42ee51fcad19cc5d0a141975e9ece5477/e2e/$guid2.md", "", "", "$guid2.md") | Out-Null

# ----------------------------------------------------------------------
# de-de sheet
# ----------------------------------------------------------------------
$hash1de = $hash1zh
$hash2de = $hash2zh

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Hyperlinks.Delete()

$loDe = $wsDe.ListObjects.Item(1)
$rowDe3 = $loDe.ListRows.Add()

# Row 2 (existing file re-run with a new GUID name)
$wsDe.Range("A2").Value = "$guid1.md"
$wsDe.Range("B2").Value = ".md"
$wsDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDe.Range("D2").Value = "e2e"
$wsDe.Range("E2").Value = "ht"
$wsDe.Range("F2").Value = "False"
$wsDe.Range("G2").Value = "$guid1.$hash1de.de-de.xlf"
$wsDe.Range("H2").Value = "2016-08-12 17:16:36"
$wsDe.Range("I2").Value = "$guid1.md"
$wsDe.Range("J2").Value = "$guid1.$hash1de.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-12 17:17:21"
$wsDe.Range("L2").Value = ""
$wsDe.Range("M2").Value = "True"
$wsDe.Range("N2").Value = ""
$wsDe.Range("O2").Value = "False"
$wsDe.Range("P2").Value = ""

# Row 3 (new file)
$wsDe.Range("A3").Value = "$guid2.md"
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "True"
$wsDe.Range("G3").Value = "$guid2.$hash2de.de-de.xlf"
$wsDe.Range("H3").Value = "2016-08-12 17:16:36"
$wsDe.Range("I3").Value = "$guid2.md"
$wsDe.Range("J3").Value = "$guid2.$hash2de.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-12 17:17:21"
$wsDe.Range("L3").Value = ""
$wsDe.Range("M3").Value = "True"
$wsDe.Range("N3").Value = ""
$wsDe.Range("O3").Value = "False"
$wsDe.Range("P3").Value = ""

$wsDe.Range("H2:H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("K2:K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/6fa018cb586824fa5ee6519a285314c1c36169a3/e2e/$guid1.md", "", "", "$guid1.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/51db68ca29367b0b8ecae2298b920256f149e2c6/e2e/$guid1.md", "", "", "$guid1.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/6fa018cb586824fa5ee6519a285314c1c36169a3/e2e/$guid2.md", "", "", "$guid2.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/51db68ca29367b0b8ecae2298b920256f149e2c6/e2e/$guid2.md", "", "", "$guid2.md") | Out-Null
